$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 312670
$ws.Range("D2").Value = 398708282
$ws.Range("C8").Value = 842
$ws.Range("D8").Value = 1240607
$ws.Range("C10").Value = 115236
$ws.Range("D10").Value = 168876720
$ws.Range("C12").Value = 57972
$ws.Range("D12").Value = 83682410
$ws.Range("C16").Value = 3901
$ws.Range("D16").Value = 5541112
$ws.Range("C20").Value = 6328
$ws.Range("D20").Value = 8833021
$ws.Range("C22").Value = 75750
$ws.Range("D22").Value = 94592520
$ws.Range("C23").Value = 53
$ws.Range("D23").Value = 69709
$ws.Range("C28").Value = 32025
$ws.Range("D28").Value = 46887378
$ws.Range("C30").Value = 11241
$ws.Range("D30").Value = 16176617
$ws.Range("C33").Value = 1538
$ws.Range("D33").Value = 2160632
$ws.Range("C35").Value = 1732
$ws.Range("D35").Value = 2442462
$ws.Range("C36").Value = 95220
$ws.Range("D36").Value = 120020647
$ws.Range("C37").Value = 65
$ws.Range("D37").Value = 75537
$ws.Range("C38").Value = 81
$ws.Range("D38").Value = 113508
$ws.Range("C42").Value = 896
$ws.Range("D42").Value = 1318685
$ws.Range("C44").Value = 43840
$ws.Range("D44").Value = 64265036
$ws.Range("C45").Value = 24
$ws.Range("D45").Value = 35950
$ws.Range("C46").Value = 8951
$ws.Range("D46").Value = 12848358
$ws.Range("C48").Value = 1378
$ws.Range("D48").Value = 1913400
$ws.Range("C51").Value = 2205
$ws.Range("D51").Value = 3073244
$ws.Range("C52").Value = 67506
$ws.Range("D52").Value = 84750976
$ws.Range("C56").Value = 380
$ws.Range("D56").Value = 558465
$ws.Range("C58").Value = 27715
$ws.Range("D58").Value = 40647682
$ws.Range("C61").Value = 10847
$ws.Range("D61").Value = 15685809
$ws.Range("C63").Value = 1338
$ws.Range("D63").Value = 1869289
$ws.Range("C67").Value = 1404
$ws.Range("D67").Value = 1963318
$ws.Range("C69").Value = 19980
$ws.Range("D69").Value = 26168565
$ws.Range("C73").Value = 7415
$ws.Range("D73").Value = 10856094
$ws.Range("C75").Value = 4973
$ws.Range("D75").Value = 7219515
$ws.Range("C76").Value = 470
$ws.Range("D76").Value = 663739
$ws.Range("C77").Value = 266
$ws.Range("D77").Value = 372673
$ws.Range("C78").Value = 137547
$ws.Range("D78").Value = 171647175
$ws.Range("C82").Value = 422
$ws.Range("D82").Value = 616416
$ws.Range("C84").Value = 62603
$ws.Range("D84").Value = 91763728
$ws.Range("C85").Value = 79
$ws.Range("D85").Value = 117082
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = 13500
$ws.Range("C87").Value = 29066
$ws.Range("D87").Value = 42061809
$ws.Range("C89").Value = 2659
$ws.Range("D89").Value = 3827740
$ws.Range("C90").Value = 2689
$ws.Range("D90").Value = 3799740
$ws.Range("C91").Value = 31158
$ws.Range("D91").Value = 42216821
$ws.Range("C94").Value = 26
$ws.Range("D94").Value = 37314
$ws.Range("C95").Value = 7675
$ws.Range("D95").Value = 11288846
$ws.Range("C97").Value = 6917
$ws.Range("D97").Value = 10027028
$ws.Range("C99").Value = 506
$ws.Range("D99").Value = 721885
$ws.Range("C100").Value = 459
$ws.Range("D100").Value = 662524
$ws.Range("C101").Value = 8310
$ws.Range("D101").Value = 11537874
$ws.Range("C103").Value = 2120
$ws.Range("D103").Value = 3122007
$ws.Range("C105").Value = 2828
$ws.Range("D105").Value = 4131362
$ws.Range("C107").Value = 116
$ws.Range("D107").Value = 168120
$ws.Range("C108").Value = 156
$ws.Range("D108").Value = 221265
$ws.Range("C109").Value = 138085
$ws.Range("D109").Value = 170829419
$ws.Range("C113").Value = 943
$ws.Range("D113").Value = 1384882
$ws.Range("C115").Value = 51919
$ws.Range("D115").Value = 76122301
$ws.Range("C117").Value = 26334
$ws.Range("D117").Value = 38153062
$ws.Range("C118").Value = 1283
$ws.Range("D118").Value = 1755065
$ws.Range("C121").Value = 2154
$ws.Range("D121").Value = 3024553
$ws.Range("C123").Value = 485191
$ws.Range("D123").Value = 639993867
$ws.Range("C125").Value = 205
$ws.Range("D125").Value = 302496
$ws.Range("C128").Value = 1349
$ws.Range("D128").Value = 1999811
$ws.Range("C130").Value = 202363
$ws.Range("D130").Value = 297524598
$ws.Range("C131").Value = 383
$ws.Range("D131").Value = 571290
$ws.Range("C133").Value = 174505
$ws.Range("D133").Value = 253662759
$ws.Range("C136").Value = 2752
$ws.Range("D136").Value = 3864704
$ws.Range("C138").Value = 6039
$ws.Range("D138").Value = 8530294
$ws.Range("C141").Value = 43177
$ws.Range("D141").Value = 57671497
$ws.Range("C146").Value = 424
$ws.Range("D146").Value = 635665
$ws.Range("C147").Value = 13751
$ws.Range("D147").Value = 20177428
$ws.Range("C148").Value = 3652
$ws.Range("D148").Value = 5268267
$ws.Range("C150").Value = 6
$ws.Range("D150").Value = 9000
$ws.Range("C151").Value = 383
$ws.Range("D151").Value = 550502
$ws.Range("C153").Value = 359
$ws.Range("D153").Value = 505510
$ws.Range("C154").Value = 16916
$ws.Range("D154").Value = 22356351
$ws.Range("C157").Value = 53
$ws.Range("D157").Value = 77906
$ws.Range("C158").Value = 6936
$ws.Range("D158").Value = 10091531
$ws.Range("C160").Value = 4818
$ws.Range("D160").Value = 6938965
$ws.Range("C162").Value = 268
$ws.Range("D162").Value = 370035
$ws.Range("C163").Value = 252
$ws.Range("D163").Value = 360783
$ws.Range("C165").Value = 14355
$ws.Range("D165").Value = 20817296
$ws.Range("C166").Value = 1689
$ws.Range("D166").Value = 2512130
$ws.Range("C167").Value = 227
$ws.Range("D167").Value = 335302
$ws.Range("C169").Value = 43
$ws.Range("D169").Value = 64190
$ws.Range("C171").Value = 85921
$ws.Range("D171").Value = 107518453
$ws.Range("C173").Value = 87
$ws.Range("D173").Value = 124454
$ws.Range("C176").Value = 636
$ws.Range("D176").Value = 937348
$ws.Range("C178").Value = 33376
$ws.Range("D178").Value = 48953841
$ws.Range("C180").Value = 12734
$ws.Range("D180").Value = 18398319
$ws.Range("C182").Value = 1218
$ws.Range("D182").Value = 1704817
$ws.Range("C184").Value = 1577
$ws.Range("D184").Value = 2219199
$ws.Range("C186").Value = 232828
$ws.Range("D186").Value = 289564744
$ws.Range("C194").Value = 85353
$ws.Range("D194").Value = 125129958
$ws.Range("C197").Value = 32332
$ws.Range("D197").Value = 46537322
$ws.Range("C200").Value = 4976
$ws.Range("D200").Value = 7091639
$ws.Range("C203").Value = 4616
$ws.Range("D203").Value = 6387299
$ws.Range("C206").Value = 256921
$ws.Range("D206").Value = 318073487
$ws.Range("C207").Value = 153
$ws.Range("D207").Value = 168018
$ws.Range("C213").Value = 605
$ws.Range("D213").Value = 881406
$ws.Range("C215").Value = 93514
$ws.Range("D215").Value = 136827854
$ws.Range("C216").Value = 85
$ws.Range("D216").Value = 126699
$ws.Range("C218").Value = 50262
$ws.Range("D218").Value = 72650140
$ws.Range("C221").Value = 4550
$ws.Range("D221").Value = 6384586
$ws.Range("C224").Value = 5435
$ws.Range("D224").Value = 7517354
$ws.Range("C227").Value = 103880
$ws.Range("D227").Value = 130097755
$ws.Range("C228").Value = 73
$ws.Range("D228").Value = 77405
$ws.Range("C234").Value = 48792
$ws.Range("D234").Value = 71487994
$ws.Range("C236").Value = 12101
$ws.Range("D236").Value = 17396509
$ws.Range("C238").Value = 1869
$ws.Range("D238").Value = 2678609
$ws.Range("C240").Value = 2401
$ws.Range("D240").Value = 3353096
$ws.Range("C241").Value = 251034
$ws.Range("D241").Value = 317146423
$ws.Range("C242").Value = 168
$ws.Range("D242").Value = 207959
$ws.Range("C244").Value = 14
$ws.Range("D244").Value = 21000
$ws.Range("C247").Value = 814
$ws.Range("D247").Value = 1195550
$ws.Range("C249").Value = 94160
$ws.Range("D249").Value = 137986579
$ws.Range("C252").Value = 63402
$ws.Range("D252").Value = 91886018
$ws.Range("C254").Value = 2350
$ws.Range("D254").Value = 3317589
$ws.Range("C257").Value = 4387
$ws.Range("D257").Value = 6156842
